$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 271.85715
$ws.Range("I28").Value = 271.85715
$ws.Range("K28").Value = 271.85715
$ws.Range("M28").Value = 213.14285

$ws.Range("H107").Value = 6142.857
$ws.Range("I107").Value = 4500
$ws.Range("J107").Value = 8333.333000000001
$ws.Range("K107").Value = 4500
$ws.Range("L107").Value = 8333.333000000001
$ws.Range("M107").Value = -2580
$ws.Range("N107").Value = -12173.333

$ws.Range("H112").Value = 2966.5557
$ws.Range("I112").Value = 1233
$ws.Range("J112").Value = 3833.3333
$ws.Range("K112").Value = 3699
$ws.Range("L112").Value = 11499.9999
$ws.Range("M112").Value = -2591
$ws.Range("N112").Value = -13715.9999

$ws.Range("H131").Value = 1775.25
$ws.Range("I131").Value = 330
$ws.Range("K131").Value = 990
$ws.Range("M131").Value = 4050

$ws.Range("H137").Value = 4494.5293
$ws.Range("I137").Value = 4163.625
$ws.Range("K137").Value = 12490.875
$ws.Range("M137").Value = -9940.875

$ws.Range("H138").Value = 2670.8333
$ws.Range("I138").Value = 2012.5
$ws.Range("K138").Value = 6037.5
$ws.Range("M138").Value = -897.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4207.8
$ws.Range("I32").Value = 1477.3438
$ws.Range("K32").Value = 1477.3438
$ws.Range("M32").Value = -1190.3438

$ws.Range("H61").Value = 3282.1428
$ws.Range("I61").Value = 2995.2
$ws.Range("K61").Value = 2995.2
$ws.Range("M61").Value = -2783.2

$ws.Range("H63").Value = 1593.25
$ws.Range("I63").Value = 1593.25
$ws.Range("K63").Value = 1593.25
$ws.Range("M63").Value = -907.25

$ws.Range("H66").Value = 1593.25
$ws.Range("I66").Value = 1593.25
$ws.Range("K66").Value = 7966.25
$ws.Range("M66").Value = -4534.25

$ws.Range("H110").Value = 2940.8333
$ws.Range("I110").Value = 3161.25
$ws.Range("J110").Value = 2500
$ws.Range("K110").Value = 3161.25
$ws.Range("L110").Value = 2500
$ws.Range("M110").Value = -1116.25
$ws.Range("N110").Value = -6590

$ws.Range("H122").Value = 3012
$ws.Range("I122").Value = 3012
$ws.Range("K122").Value = 9036
$ws.Range("M122").Value = -6586

$ws.Range("H132").Value = 2431.1428
$ws.Range("I132").Value = 2419.6667
$ws.Range("K132").Value = 7259.000100000001
$ws.Range("M132").Value = -4729.000100000001

$ws.Range("H136").Value = 3282.1428
$ws.Range("I136").Value = 2995.2
$ws.Range("K136").Value = 8985.599999999999
$ws.Range("M136").Value = -6435.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H80").Value = 145.61539
$ws.Range("I80").Value = 75
$ws.Range("K80").Value = 75
$ws.Range("M80").Value = 923

$ws.Range("H83").Value = 145.61539
$ws.Range("I83").Value = 75
$ws.Range("K83").Value = 375
$ws.Range("M83").Value = 4617

$ws.Range("H103").Value = 4800
$ws.Range("J103").Value = 4800
$ws.Range("L103").Value = 4800
$ws.Range("N103").Value = -7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 375
$ws.Range("I22").Value = 333.33334
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 333.33334
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 16.66665999999998
$ws.Range("N22").Value = -1200

$ws.Range("H99").Value = 1804
$ws.Range("I99").Value = 1456
$ws.Range("K99").Value = 1456
$ws.Range("M99").Value = 42

$ws.Range("H122").Value = 2398.75
$ws.Range("I122").Value = 2112.3333
$ws.Range("K122").Value = 6336.999899999999
$ws.Range("M122").Value = -3886.999899999999

$ws.Range("H126").Value = 1804
$ws.Range("I126").Value = 1456
$ws.Range("K126").Value = 4368
$ws.Range("M126").Value = -1898

$ws.Range("H132").Value = 2682.4
$ws.Range("I132").Value = 699
$ws.Range("J132").Value = 4004.6667
$ws.Range("K132").Value = 2097
$ws.Range("L132").Value = 12014.0001
$ws.Range("M132").Value = 433
$ws.Range("N132").Value = -17074.0001

$ws.Range("H134").Value = 6252
$ws.Range("I134").Value = 6252
$ws.Range("K134").Value = 18756
$ws.Range("M134").Value = -16221

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 333417
$ws.Range("I11").Value = 500050
$ws.Range("K11").Value = 1500150
$ws.Range("M11").Value = -1500010

$ws.Range("H16").Value = 519.5
$ws.Range("I16").Value = 519.5
$ws.Range("K16").Value = 1558.5
$ws.Range("M16").Value = -1385.5

$ws.Range("H19").Value = 361.66666
$ws.Range("I19").Value = 10
$ws.Range("J19").Value = 537.5
$ws.Range("K19").Value = 30
$ws.Range("L19").Value = 1612.5
$ws.Range("M19").Value = 144
$ws.Range("N19").Value = -1960.5

$ws.Range("H34").Value = 1863.0834
$ws.Range("I34").Value = 262.2
$ws.Range("J34").Value = 3006.5715
$ws.Range("K34").Value = 786.5999999999999
$ws.Range("L34").Value = 9019.7145
$ws.Range("M34").Value = -702.5999999999999
$ws.Range("N34").Value = -9187.7145

$ws.Range("H74").Value = 258331.17
$ws.Range("J74").Value = 303994.8
$ws.Range("L74").Value = 911984.3999999999
$ws.Range("N74").Value = -914106.3999999999

$ws.Range("H77").Value = 258331.17
$ws.Range("J77").Value = 303994.8
$ws.Range("L77").Value = 2735953.2
$ws.Range("N77").Value = -2746561.2

$ws.Range("H97").Value = 404.6
$ws.Range("J97").Value = 752
$ws.Range("L97").Value = 2256
$ws.Range("N97").Value = -3248

$ws.Range("H122").Value = 995
$ws.Range("I122").Value = 995
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8955
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6505
$ws.Range("N122").ClearContents()

$ws.Range("H128").Value = 308165.5
$ws.Range("I128").Value = 308165.5
$ws.Range("K128").Value = 924496.5
$ws.Range("M128").Value = -919516.5

$ws.Range("H131").Value = 1383.7778
$ws.Range("I131").Value = 993.5
$ws.Range("J131").Value = 1696
$ws.Range("K131").Value = 2980.5
$ws.Range("L131").Value = 5088
$ws.Range("M131").Value = 2059.5
$ws.Range("N131").Value = -15168

$ws.Range("H141").Value = 2750
$ws.Range("I141").Value = 2750
$ws.Range("K141").Value = 8250
$ws.Range("M141").Value = -3070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2291.5293
$ws.Range("I102").Value = 2243.4666
$ws.Range("K102").Value = 2243.4666
$ws.Range("M102").Value = -621.4666000000002

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5209.8
$ws.Range("I40").Value = 4514.6665
$ws.Range("K40").Value = 4514.6665
$ws.Range("M40").Value = -4378.6665

$ws.Range("H46").Value = 3312.25
$ws.Range("I46").Value = 2789.8
$ws.Range("J46").Value = 3685.4285
$ws.Range("K46").Value = 2789.8
$ws.Range("L46").Value = 3685.4285
$ws.Range("M46").Value = -2601.8
$ws.Range("N46").Value = -4061.4285

$ws.Range("H82").Value = 1024.2
$ws.Range("I82").Value = 1024.2
$ws.Range("K82").Value = 1024.2
$ws.Range("M82").Value = -663.2

$ws.Range("H85").Value = 1024.2
$ws.Range("I85").Value = 1024.2
$ws.Range("K85").Value = 1024.2
$ws.Range("M85").Value = 223.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3700
$ws.Range("I62").Value = 3400
$ws.Range("K62").Value = 3400
$ws.Range("M62").Value = -2776

$ws.Range("H65").Value = 3700
$ws.Range("I65").Value = 3400
$ws.Range("K65").Value = 17000
$ws.Range("M65").Value = -13880

$ws.Range("H122").Value = 1561.8
$ws.Range("I122").Value = 1104.6666
$ws.Range("K122").Value = 3313.9998
$ws.Range("M122").Value = -863.9998000000001

$ws.Range("H126").Value = 39411
$ws.Range("I126").Value = 38671.215
$ws.Range("K126").Value = 116013.645
$ws.Range("M126").Value = -113543.645

$ws.Range("H136").Value = 10899.5
$ws.Range("I136").Value = 10899.5
$ws.Range("K136").Value = 32698.5
$ws.Range("M136").Value = -30148.5
